$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.293.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.02%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.928.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7568'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.42%  '

# Row 7
$ws.Range("E7").Value = '  -0.21%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.15'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.62%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3179'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07033'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08021'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.933.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.376'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.34%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.46%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.77%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.290.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '253.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.27%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.884'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.46%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007969'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.183.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.705'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.64%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.502'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.58%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1337'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.213'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.368'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.84%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.512'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.411'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.135'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.23%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05233'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.14%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.317'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7544'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.42%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.787'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01957'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.796'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.41%  '

# Row 40
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.24%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.520'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.48%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4489'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.972'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.95%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.20%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8348'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.84%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.934'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.91%  '

# Row 47
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.600'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.20%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.55%  '

# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '982.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1207'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.53%  '
